# SWE Group 4 Charter Document - apply "did 4, 5, signed / My part" edit.
#
# The charter template (numbered ListParagraph bullets 1-8 plus signature
# block) is replaced by the author's actual filled-in answers for items
# 4-6 plus a single signature ("Cayden Hannon"). The new paragraphs carry
# pasted-in direct formatting (Arial/Lato, sz 25, color 000000, kern 0,
# ligatures disabled) instead of the ListParagraph/numbered-list style.
#
# We rebuild the whole document body in one shot via Range.InsertXML,
# which replaces the target range's contents with the supplied WordML -
# letting us set the exact run-level rPr (including the w14:ligatures
# extension element, which isn't reachable through Font.Ligatures).

$d = $word.ActiveDocument

$bodyXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
  <w:pPr>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="000000"/>
      <w:kern w:val="0"/>
      <w:sz w:val="25"/>
      <w:szCs w:val="25"/>
      <w14:ligatures w14:val="none"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="000000"/>
      <w:kern w:val="0"/>
      <w:sz w:val="25"/>
      <w:szCs w:val="25"/>
      <w14:ligatures w14:val="none"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Lato" w:eastAsia="Times New Roman" w:hAnsi="Lato" w:cs="Times New Roman"/>
      <w:color w:val="000000"/>
      <w:kern w:val="0"/>
      <w14:ligatures w14:val="none"/>
    </w:rPr>
    <w:br/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="000000"/>
      <w:kern w:val="0"/>
      <w:sz w:val="25"/>
      <w:szCs w:val="25"/>
      <w14:ligatures w14:val="none"/>
    </w:rPr>
    <w:t>5. Procedure for adopting policies (</w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="000000"/>
      <w:kern w:val="0"/>
      <w:sz w:val="25"/>
      <w:szCs w:val="25"/>
      <w14:ligatures w14:val="none"/>
    </w:rPr>
    <w:t>e.g.</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="000000"/>
      <w:kern w:val="0"/>
      <w:sz w:val="25"/>
      <w:szCs w:val="25"/>
      <w14:ligatures w14:val="none"/>
    </w:rPr>
    <w:t xml:space="preserve"> code style, code check-in steps, documentation</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Lato" w:eastAsia="Times New Roman" w:hAnsi="Lato" w:cs="Times New Roman"/>
      <w:color w:val="000000"/>
      <w:kern w:val="0"/>
      <w14:ligatures w14:val="none"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="000000"/>
      <w:kern w:val="0"/>
      <w:sz w:val="25"/>
      <w:szCs w:val="25"/>
      <w14:ligatures w14:val="none"/>
    </w:rPr>
    <w:t>updates)</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="000000"/>
      <w:kern w:val="0"/>
      <w:sz w:val="25"/>
      <w:szCs w:val="25"/>
      <w14:ligatures w14:val="none"/>
    </w:rPr>
    <w:t xml:space="preserve">: </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="000000"/>
      <w:kern w:val="0"/>
      <w:sz w:val="25"/>
      <w:szCs w:val="25"/>
      <w14:ligatures w14:val="none"/>
    </w:rPr>
    <w:t>Group vote online, ensure that everyone is informed of the policy. Also make sure it is possible by every member.</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Lato" w:eastAsia="Times New Roman" w:hAnsi="Lato" w:cs="Times New Roman"/>
      <w:color w:val="000000"/>
      <w:kern w:val="0"/>
      <w14:ligatures w14:val="none"/>
    </w:rPr>
    <w:br/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="000000"/>
      <w:kern w:val="0"/>
      <w:sz w:val="25"/>
      <w:szCs w:val="25"/>
      <w14:ligatures w14:val="none"/>
    </w:rPr>
    <w:t>6. if you adopt such a policy during the initial meeting, put it here</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="000000"/>
      <w:kern w:val="0"/>
      <w:sz w:val="25"/>
      <w:szCs w:val="25"/>
      <w14:ligatures w14:val="none"/>
    </w:rPr>
    <w:t xml:space="preserve">: </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="000000"/>
      <w:kern w:val="0"/>
      <w:sz w:val="25"/>
      <w:szCs w:val="25"/>
      <w14:ligatures w14:val="none"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="000000"/>
      <w:kern w:val="0"/>
      <w:sz w:val="25"/>
      <w:szCs w:val="25"/>
      <w14:ligatures w14:val="none"/>
    </w:rPr>
    <w:t>N/A</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="000000"/>
      <w:kern w:val="0"/>
      <w:sz w:val="25"/>
      <w:szCs w:val="25"/>
      <w14:ligatures w14:val="none"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:ascii="Lato" w:eastAsia="Times New Roman" w:hAnsi="Lato" w:cs="Times New Roman"/>
      <w:color w:val="000000"/>
      <w:kern w:val="0"/>
      <w14:ligatures w14:val="none"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:kern w:val="0"/>
      <w14:ligatures w14:val="none"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:r>
    <w:t>Cayden Hannon</w:t>
  </w:r>
</w:p>
'@

# Replacing $d.Content leaves the sectPr (page setup/header ref) alone and
# collapses every existing paragraph down to just the ones supplied here.
$d.Content.InsertXML($bodyXml)

# styles.xml also picks up an (unused) character style in the real commit,
# presumably left behind by pasted-in web content; add it for fidelity.
$newStyle = $d.Styles.Add("textlayer--absolute", 2)
$newStyle.BaseStyle = $d.Styles("DefaultParagraphFont")

Write-Host ("Paragraphs: " + $d.Paragraphs.Count)
